# Update cryptocurrency price/volume data in cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to keep a plain text representation (matches the
    # original inline-string cells) instead of letting Excel auto-convert
    # numeric-looking strings (e.g. "0.9997") into real numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.524.17"
$ws.Range("E2").Value = "  -0.33%  "

Set-TextValue $ws.Range("D3") "1.910.85"
$ws.Range("E3").Value = "  -0.64%  "

Set-TextValue $ws.Range("D4") "0.9997"
$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue $ws.Range("D5") "244.12"
$ws.Range("E5").Value = "  -1.22%  "

Set-TextValue $ws.Range("D6") "0.9995"
$ws.Range("E6").Value = "  -0.06%  "

Set-TextValue $ws.Range("D7") "0.4844"
$ws.Range("E7").Value = "  +2.16%  "

Set-TextValue $ws.Range("D8") "0.2893"
$ws.Range("E8").Value = "  +0.14%  "

Set-TextValue $ws.Range("D9") "0.06801"
$ws.Range("E9").Value = "  -0.68%  "

Set-TextValue $ws.Range("D10") "111.08"
$ws.Range("E10").Value = "  +5.51%  "

Set-TextValue $ws.Range("D11") "19.27"
$ws.Range("E11").Value = "  +4.87%  "

Set-TextValue $ws.Range("D12") "1.912.52"
$ws.Range("E12").Value = "  -0.53%  "

Set-TextValue $ws.Range("D13") "0.07564"
$ws.Range("E13").Value = "  -1.75%  "

Set-TextValue $ws.Range("D14") "5.378"
$ws.Range("E14").Value = "  +0.92%  "

Set-TextValue $ws.Range("D15") "0.6704"
$ws.Range("E15").Value = "  +0.27%  "

Set-TextValue $ws.Range("D16") "296.68"
$ws.Range("E16").Value = "  +1.60%  "

Set-TextValue $ws.Range("D17") "30.528.85"

$ws.Range("E18").Value = "  +0.45%  "

Set-TextValue $ws.Range("D19") "0.000007589"
$ws.Range("E19").Value = "  -0.43%  "

Set-TextValue $ws.Range("D20") "0.9994"
$ws.Range("E20").Value = "  -0.05%  "

Set-TextValue $ws.Range("D21") "5.522"
$ws.Range("E21").Value = "  -1.22%  "

Set-TextValue $ws.Range("D22") "2.160.94"
$ws.Range("E22").Value = "  -0.65%  "

Set-TextValue $ws.Range("D23") "1.001"
$ws.Range("E23").Value = "  +0.06%  "

Set-TextValue $ws.Range("D24") "6.431"
$ws.Range("E24").Value = "  -0.38%  "

Set-TextValue $ws.Range("D25") "9.461"
$ws.Range("E25").Value = "  -0.07%  "

Set-TextValue $ws.Range("D26") "165.70"
$ws.Range("E26").Value = "  -1.17%  "

Set-TextValue $ws.Range("D27") "20.32"
$ws.Range("E27").Value = "  -3.54%  "

Set-TextValue $ws.Range("D28") "2.080"
$ws.Range("E28").Value = "  -1.95%  "

$ws.Range("E29").Value = "  -0.65%  "

Set-TextValue $ws.Range("D30") "1.434"
$ws.Range("E30").Value = "  +2.90%  "

Set-TextValue $ws.Range("D31") "4.145"

Set-TextValue $ws.Range("D32") "4.043"
$ws.Range("E32").Value = "  -0.59%  "

Set-TextValue $ws.Range("D34") "0.7363"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("E35").Value = "  -0.78%  "

Set-TextValue $ws.Range("D36") "0.9993"
$ws.Range("E36").Value = "  +0.02%  "

Set-TextValue $ws.Range("D37") "0.02036"
$ws.Range("E37").Value = "  -1.64%  "

Set-TextValue $ws.Range("D38") "2.713"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("E39").Value = "  -0.37%  "

Set-TextValue $ws.Range("D40") "2.020"
$ws.Range("E40").Value = "  -1.81%  "

Set-TextValue $ws.Range("D41") "109.20"
$ws.Range("E41").Value = "  -1.90%  "

Set-TextValue $ws.Range("D42") "0.4451"
$ws.Range("E42").Value = "  +1.53%  "

Set-TextValue $ws.Range("D43") "0.8674"
$ws.Range("E43").Value = "  -1.07%  "

Set-TextValue $ws.Range("D44") "5.785"
$ws.Range("E44").Value = "  -2.10%  "

Set-TextValue $ws.Range("D45") "0.9991"
$ws.Range("E45").Value = "  -0.09%  "

Set-TextValue $ws.Range("D46") "69.38"
$ws.Range("E46").Value = "  +1.92%  "

Set-TextValue $ws.Range("D47") "7.196"
$ws.Range("E47").Value = "  -1.24%  "

Set-TextValue $ws.Range("D48") "48.35"
$ws.Range("E48").Value = "  +0.24%  "

Set-TextValue $ws.Range("D49") "9.235"
$ws.Range("E49").Value = "  -1.58%  "

Set-TextValue $ws.Range("D50") "0.1226"
$ws.Range("E50").Value = "  -1.43%  "

Set-TextValue $ws.Range("D51") "0.2509"
$ws.Range("E51").Value = "  -0.79%  "
